# Auto-generated edit script: updates Louisoix_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the commit diff. Values are set via the Excel COM object model; cells that became
# empty in the target are cleared with ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 695.7377300000001  # ALC!H17: 695.7213 -> 695.7377300000001
$ws.Cells.Item(17, 9).Value = 288  # ALC!I17: 0 -> 288
$ws.Cells.Item(17, 10).Value = 702.5333000000001  # ALC!J17: 695.7213 -> 702.5333000000001
$ws.Cells.Item(17, 11).Value = 864  # ALC!K17: 0 -> 864
$ws.Cells.Item(17, 12).Value = 2107.5999  # ALC!L17: 2087.1639 -> 2107.5999
$ws.Cells.Item(17, 13).Value = -696  # ALC!M17: None -> -696
$ws.Cells.Item(17, 14).Value = -2443.5999  # ALC!N17: -2423.1639 -> -2443.5999

$ws.Cells.Item(20, 8).Value = 1007  # ALC!H20: 0 -> 1007
$ws.Cells.Item(20, 9).Value = 1007  # ALC!I20: 0 -> 1007
$ws.Cells.Item(20, 11).Value = 1007  # ALC!K20: 0 -> 1007
$ws.Cells.Item(20, 13).Value = -777  # ALC!M20: None -> -777

$ws.Cells.Item(28, 8).Value = 2751.3333  # ALC!H28: 758.125 -> 2751.3333
$ws.Cells.Item(28, 9).Value = 680.7143  # ALC!I28: 758.125 -> 680.7143
$ws.Cells.Item(28, 10).Value = 9998.5  # ALC!J28: 0 -> 9998.5
$ws.Cells.Item(28, 11).Value = 680.7143  # ALC!K28: 758.125 -> 680.7143
$ws.Cells.Item(28, 12).Value = 9998.5  # ALC!L28: 0 -> 9998.5
$ws.Cells.Item(28, 13).Value = -195.7143  # ALC!M28: -273.125 -> -195.7143
$ws.Cells.Item(28, 14).Value = -10968.5  # ALC!N28: None -> -10968.5

$ws.Cells.Item(35, 8).Value = 1007  # ALC!H35: 0 -> 1007
$ws.Cells.Item(35, 9).Value = 1007  # ALC!I35: 0 -> 1007
$ws.Cells.Item(35, 11).Value = 1007  # ALC!K35: 0 -> 1007
$ws.Cells.Item(35, 13).Value = -628  # ALC!M35: None -> -628

$ws.Cells.Item(40, 8).Value = 2636.923  # ALC!H40: 2797.875 -> 2636.923
$ws.Cells.Item(40, 9).Value = 1374.1111  # ALC!I40: 1436.6875 -> 1374.1111
$ws.Cells.Item(40, 10).Value = 5478.25  # ALC!J40: 5520.25 -> 5478.25
$ws.Cells.Item(40, 11).Value = 1374.1111  # ALC!K40: 1436.6875 -> 1374.1111
$ws.Cells.Item(40, 12).Value = 5478.25  # ALC!L40: 5520.25 -> 5478.25
$ws.Cells.Item(40, 13).Value = -1199.1111  # ALC!M40: -1261.6875 -> -1199.1111
$ws.Cells.Item(40, 14).Value = -5828.25  # ALC!N40: -5870.25 -> -5828.25

$ws.Cells.Item(53, 8).Value = 517.7857  # ALC!H53: 443.2857 -> 517.7857
$ws.Cells.Item(53, 9).Value = 157  # ALC!I53: 146.90909 -> 157
$ws.Cells.Item(53, 10).Value = 1419.75  # ALC!J53: 1530 -> 1419.75
$ws.Cells.Item(53, 11).Value = 157  # ALC!K53: 146.90909 -> 157
$ws.Cells.Item(53, 12).Value = 1419.75  # ALC!L53: 1530 -> 1419.75
$ws.Cells.Item(53, 13).Value = 480  # ALC!M53: 490.09091 -> 480
$ws.Cells.Item(53, 14).Value = -2693.75  # ALC!N53: -2804 -> -2693.75

$ws.Cells.Item(55, 8).Value = 501.1  # ALC!H55: 553.05554 -> 501.1
$ws.Cells.Item(55, 9).Value = 229.25  # ALC!I55: 294.5 -> 229.25
$ws.Cells.Item(55, 11).Value = 229.25  # ALC!K55: 294.5 -> 229.25
$ws.Cells.Item(55, 13).Value = -15.25  # ALC!M55: -80.5 -> -15.25

$ws.Cells.Item(70, 8).Value = 2993.5  # ALC!H70: 2810.75 -> 2993.5
$ws.Cells.Item(70, 10).Value = 2993.5  # ALC!J70: 2810.75 -> 2993.5
$ws.Cells.Item(70, 12).Value = 8980.5  # ALC!L70: 8432.25 -> 8980.5
$ws.Cells.Item(70, 14).Value = -9520.5  # ALC!N70: -8972.25 -> -9520.5

$ws.Cells.Item(73, 8).Value = 2993.5  # ALC!H73: 2810.75 -> 2993.5
$ws.Cells.Item(73, 10).Value = 2993.5  # ALC!J73: 2810.75 -> 2993.5
$ws.Cells.Item(73, 12).Value = 8980.5  # ALC!L73: 8432.25 -> 8980.5
$ws.Cells.Item(73, 14).Value = -10852.5  # ALC!N73: -10304.25 -> -10852.5

$ws.Cells.Item(76, 8).Value = 6751  # ALC!H76: 5180.8184 -> 6751
$ws.Cells.Item(76, 9).Value = 5000  # ALC!I76: 4373.125 -> 5000
$ws.Cells.Item(76, 11).Value = 5000  # ALC!K76: 4373.125 -> 5000
$ws.Cells.Item(76, 13).Value = -4685  # ALC!M76: -4058.125 -> -4685

$ws.Cells.Item(79, 8).Value = 6751  # ALC!H79: 5180.8184 -> 6751
$ws.Cells.Item(79, 9).Value = 5000  # ALC!I79: 4373.125 -> 5000
$ws.Cells.Item(79, 11).Value = 5000  # ALC!K79: 4373.125 -> 5000
$ws.Cells.Item(79, 13).Value = -3908  # ALC!M79: -3281.125 -> -3908

$ws.Cells.Item(86, 8).Value = 4579.769  # ALC!H86: 3796.4736 -> 4579.769
$ws.Cells.Item(86, 9).Value = 4899.25  # ALC!I86: 3699.2856 -> 4899.25
$ws.Cells.Item(86, 11).Value = 4899.25  # ALC!K86: 3699.2856 -> 4899.25
$ws.Cells.Item(86, 13).Value = -3776.25  # ALC!M86: -2576.2856 -> -3776.25

$ws.Cells.Item(89, 8).Value = 4579.769  # ALC!H89: 3796.4736 -> 4579.769
$ws.Cells.Item(89, 9).Value = 4899.25  # ALC!I89: 3699.2856 -> 4899.25
$ws.Cells.Item(89, 11).Value = 24496.25  # ALC!K89: 18496.428 -> 24496.25
$ws.Cells.Item(89, 13).Value = -18880.25  # ALC!M89: -12880.428 -> -18880.25

$ws.Cells.Item(98, 8).Value = 2033  # ALC!H98: 2032.1515 -> 2033
$ws.Cells.Item(98, 9).Value = 2084.2666  # ALC!I98: 2085.0688 -> 2084.2666
$ws.Cells.Item(98, 11).Value = 2084.2666  # ALC!K98: 2085.0688 -> 2084.2666
$ws.Cells.Item(98, 13).Value = -586.2665999999999  # ALC!M98: -587.0688 -> -586.2665999999999

$ws.Cells.Item(106, 8).Value = 11127.429  # ALC!H106: 8076.8184 -> 11127.429
$ws.Cells.Item(106, 9).Value = 7746  # ALC!I106: 5242.125 -> 7746
$ws.Cells.Item(106, 11).Value = 7746  # ALC!K106: 5242.125 -> 7746
$ws.Cells.Item(106, 13).Value = -7115  # ALC!M106: -4611.125 -> -7115

$ws.Cells.Item(107, 8).Value = 1128.2  # ALC!H107: 1066.4546 -> 1128.2
$ws.Cells.Item(107, 10).Value = 4000  # ALC!J107: 2224.5 -> 4000
$ws.Cells.Item(107, 12).Value = 4000  # ALC!L107: 2224.5 -> 4000
$ws.Cells.Item(107, 14).Value = -7840  # ALC!N107: -6064.5 -> -7840

$ws.Cells.Item(112, 8).Value = 6398.8  # ALC!H112: 6399.2 -> 6398.8
$ws.Cells.Item(112, 10).Value = 6331.3335  # ALC!J112: 6332 -> 6331.3335
$ws.Cells.Item(112, 12).Value = 18994.0005  # ALC!L112: 18996 -> 18994.0005
$ws.Cells.Item(112, 14).Value = -21210.0005  # ALC!N112: -21212 -> -21210.0005

$ws.Cells.Item(113, 8).Value = 11552.333  # ALC!H113: 7000.4 -> 11552.333
$ws.Cells.Item(113, 9).Value = 13328.833  # ALC!I113: 6667.6665 -> 13328.833
$ws.Cells.Item(113, 10).Value = 7999.3335  # ALC!J113: 7499.5 -> 7999.3335
$ws.Cells.Item(113, 11).Value = 13328.833  # ALC!K113: 6667.6665 -> 13328.833
$ws.Cells.Item(113, 12).Value = 7999.3335  # ALC!L113: 7499.5 -> 7999.3335
$ws.Cells.Item(113, 13).Value = -10074.833  # ALC!M113: -3413.6665 -> -10074.833
$ws.Cells.Item(113, 14).Value = -14507.3335  # ALC!N113: -14007.5 -> -14507.3335

$ws.Cells.Item(116, 8).Value = 16137.308  # ALC!H116: 15731.346 -> 16137.308
$ws.Cells.Item(116, 9).Value = 16597.572  # ALC!I116: 15078.375 -> 16597.572
$ws.Cells.Item(116, 10).Value = 15967.737  # ALC!J116: 16021.556 -> 15967.737
$ws.Cells.Item(116, 11).Value = 16597.572  # ALC!K116: 15078.375 -> 16597.572
$ws.Cells.Item(116, 12).Value = 15967.737  # ALC!L116: 16021.556 -> 15967.737
$ws.Cells.Item(116, 13).Value = -13155.572  # ALC!M116: -11636.375 -> -13155.572
$ws.Cells.Item(116, 14).Value = -22851.737  # ALC!N116: -22905.556 -> -22851.737

$ws.Cells.Item(122, 8).Value = 2033  # ALC!H122: 2032.1515 -> 2033
$ws.Cells.Item(122, 9).Value = 2084.2666  # ALC!I122: 2085.0688 -> 2084.2666
$ws.Cells.Item(122, 11).Value = 6252.7998  # ALC!K122: 6255.2064 -> 6252.7998
$ws.Cells.Item(122, 13).Value = -3802.7998  # ALC!M122: -3805.2064 -> -3802.7998

$ws.Cells.Item(137, 8).Value = 2872.9375  # ALC!H137: 2766.7812 -> 2872.9375
$ws.Cells.Item(137, 9).Value = 2424.45  # ALC!I137: 2289.1052 -> 2424.45
$ws.Cells.Item(137, 10).Value = 3620.4167  # ALC!J137: 3464.923 -> 3620.4167
$ws.Cells.Item(137, 11).Value = 7273.349999999999  # ALC!K137: 6867.3156 -> 7273.349999999999
$ws.Cells.Item(137, 12).Value = 10861.2501  # ALC!L137: 10394.769 -> 10861.2501
$ws.Cells.Item(137, 13).Value = -4723.349999999999  # ALC!M137: -4317.3156 -> -4723.349999999999
$ws.Cells.Item(137, 14).Value = -15961.2501  # ALC!N137: -15494.769 -> -15961.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13197.979  # ARM!H32: 12155.06 -> 13197.979
$ws.Cells.Item(32, 9).Value = 13450.823  # ARM!I32: 12307.656 -> 13450.823
$ws.Cells.Item(32, 11).Value = 13450.823  # ARM!K32: 12307.656 -> 13450.823
$ws.Cells.Item(32, 13).Value = -13163.823  # ARM!M32: -12020.656 -> -13163.823

$ws.Cells.Item(61, 8).Value = 13408.462  # ARM!H61: 11950.066 -> 13408.462
$ws.Cells.Item(61, 9).Value = 18775.125  # ARM!I61: 15454.5 -> 18775.125
$ws.Cells.Item(61, 10).Value = 4821.8  # ARM!J61: 4941.2 -> 4821.8
$ws.Cells.Item(61, 11).Value = 18775.125  # ARM!K61: 15454.5 -> 18775.125
$ws.Cells.Item(61, 12).Value = 4821.8  # ARM!L61: 4941.2 -> 4821.8
$ws.Cells.Item(61, 13).Value = -18563.125  # ARM!M61: -15242.5 -> -18563.125
$ws.Cells.Item(61, 14).Value = -5245.8  # ARM!N61: -5365.2 -> -5245.8

$ws.Cells.Item(74, 8).Value = 992.53845  # ARM!H74: 899.8 -> 992.53845
$ws.Cells.Item(74, 9).Value = 966.9524  # ARM!I74: 859.76 -> 966.9524
$ws.Cells.Item(74, 11).Value = 966.9524  # ARM!K74: 859.76 -> 966.9524
$ws.Cells.Item(74, 13).Value = -92.95240000000001  # ARM!M74: 14.24000000000001 -> -92.95240000000001

$ws.Cells.Item(77, 8).Value = 992.53845  # ARM!H77: 899.8 -> 992.53845
$ws.Cells.Item(77, 9).Value = 966.9524  # ARM!I77: 859.76 -> 966.9524
$ws.Cells.Item(77, 11).Value = 4834.762  # ARM!K77: 4298.8 -> 4834.762
$ws.Cells.Item(77, 13).Value = -466.7619999999997  # ARM!M77: 69.19999999999982 -> -466.7619999999997

$ws.Cells.Item(97, 8).Value = 2381.842  # ARM!H97: 2303.2 -> 2381.842
$ws.Cells.Item(97, 9).Value = 1211.4546  # ARM!I97: 1177.9166 -> 1211.4546
$ws.Cells.Item(97, 11).Value = 1211.4546  # ARM!K97: 1177.9166 -> 1211.4546
$ws.Cells.Item(97, 13).Value = -715.4546  # ARM!M97: -681.9166 -> -715.4546

$ws.Cells.Item(102, 8).Value = 6088.4  # ARM!H102: 6487.3335 -> 6088.4
$ws.Cells.Item(102, 9).Value = 6431.778  # ARM!I102: 6923.5 -> 6431.778
$ws.Cells.Item(102, 11).Value = 6431.778  # ARM!K102: 6923.5 -> 6431.778
$ws.Cells.Item(102, 13).Value = -4809.778  # ARM!M102: -5301.5 -> -4809.778

$ws.Cells.Item(109, 8).Value = 47038  # ARM!H109: 48110.6 -> 47038
$ws.Cells.Item(109, 10).Value = 47038  # ARM!J109: 48110.6 -> 47038
$ws.Cells.Item(109, 12).Value = 47038  # ARM!L109: 48110.6 -> 47038
$ws.Cells.Item(109, 14).Value = -49812  # ARM!N109: -50884.6 -> -49812

$ws.Cells.Item(132, 8).Value = 38877.465  # ARM!H132: 37602.31 -> 38877.465
$ws.Cells.Item(132, 9).Value = 55390.633  # ARM!I132: 52701.05 -> 55390.633
$ws.Cells.Item(132, 10).Value = 4016.3333  # ARM!J132: 4049.5557 -> 4016.3333
$ws.Cells.Item(132, 11).Value = 166171.899  # ARM!K132: 158103.15 -> 166171.899
$ws.Cells.Item(132, 12).Value = 12048.9999  # ARM!L132: 12148.6671 -> 12048.9999
$ws.Cells.Item(132, 13).Value = -163641.899  # ARM!M132: -155573.15 -> -163641.899
$ws.Cells.Item(132, 14).Value = -17108.9999  # ARM!N132: -17208.6671 -> -17108.9999

$ws.Cells.Item(133, 8).Value = 60000  # ARM!H133: 59998.332 -> 60000
$ws.Cells.Item(133, 10).Value = 60000  # ARM!J133: 59998.332 -> 60000
$ws.Cells.Item(133, 12).Value = 60000  # ARM!L133: 59998.332 -> 60000
$ws.Cells.Item(133, 14).Value = -65060  # ARM!N133: -65058.332 -> -65060

$ws.Cells.Item(136, 8).Value = 13408.462  # ARM!H136: 11950.066 -> 13408.462
$ws.Cells.Item(136, 9).Value = 18775.125  # ARM!I136: 15454.5 -> 18775.125
$ws.Cells.Item(136, 10).Value = 4821.8  # ARM!J136: 4941.2 -> 4821.8
$ws.Cells.Item(136, 11).Value = 56325.375  # ARM!K136: 46363.5 -> 56325.375
$ws.Cells.Item(136, 12).Value = 14465.4  # ARM!L136: 14823.6 -> 14465.4
$ws.Cells.Item(136, 13).Value = -53775.375  # ARM!M136: -43813.5 -> -53775.375
$ws.Cells.Item(136, 14).Value = -19565.4  # ARM!N136: -19923.6 -> -19565.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 39498.5  # BSM!H26: 49998 -> 39498.5
$ws.Cells.Item(26, 9).Value = 39498.5  # BSM!I26: 49998 -> 39498.5
$ws.Cells.Item(26, 11).Value = 39498.5  # BSM!K26: 49998 -> 39498.5
$ws.Cells.Item(26, 13).Value = -39206.5  # BSM!M26: -49706 -> -39206.5

$ws.Cells.Item(94, 8).Value = 3980.1428  # BSM!H94: 4023.6428 -> 3980.1428
$ws.Cells.Item(94, 9).Value = 1978.7693  # BSM!I94: 2025.6154 -> 1978.7693
$ws.Cells.Item(94, 11).Value = 1978.7693  # BSM!K94: 2025.6154 -> 1978.7693
$ws.Cells.Item(94, 13).Value = -1527.7693  # BSM!M94: -1574.6154 -> -1527.7693

$ws.Cells.Item(99, 8).Value = 3944.0557  # BSM!H99: 3882.4119 -> 3944.0557
$ws.Cells.Item(99, 9).Value = 2940  # BSM!I99: 2849.9 -> 2940
$ws.Cells.Item(99, 10).Value = 5952.1665  # BSM!J99: 5357.4287 -> 5952.1665
$ws.Cells.Item(99, 11).Value = 2940  # BSM!K99: 2849.9 -> 2940
$ws.Cells.Item(99, 12).Value = 5952.1665  # BSM!L99: 5357.4287 -> 5952.1665
$ws.Cells.Item(99, 13).Value = -1442  # BSM!M99: -1351.9 -> -1442
$ws.Cells.Item(99, 14).Value = -8948.166499999999  # BSM!N99: -8353.4287 -> -8948.166499999999

$ws.Cells.Item(102, 8).Value = 9859.333000000001  # BSM!H102: 8096 -> 9859.333000000001
$ws.Cells.Item(102, 9).Value = 9859.333000000001  # BSM!I102: 8096 -> 9859.333000000001
$ws.Cells.Item(102, 11).Value = 9859.333000000001  # BSM!K102: 8096 -> 9859.333000000001
$ws.Cells.Item(102, 13).Value = -6614.333000000001  # BSM!M102: -4851 -> -6614.333000000001

$ws.Cells.Item(105, 8).Value = 3509.7742  # BSM!H105: 3540.5334 -> 3509.7742
$ws.Cells.Item(105, 9).Value = 3346.9614  # BSM!I105: 3377.16 -> 3346.9614
$ws.Cells.Item(105, 10).Value = 4356.4  # BSM!J105: 4357.4 -> 4356.4
$ws.Cells.Item(105, 11).Value = 3346.9614  # BSM!K105: 3377.16 -> 3346.9614
$ws.Cells.Item(105, 12).Value = 4356.4  # BSM!L105: 4357.4 -> 4356.4
$ws.Cells.Item(105, 13).Value = -1599.9614  # BSM!M105: -1630.16 -> -1599.9614
$ws.Cells.Item(105, 14).Value = -7850.4  # BSM!N105: -7851.4 -> -7850.4

$ws.Cells.Item(134, 8).Value = 3484.862  # BSM!H134: 3631.926 -> 3484.862
$ws.Cells.Item(134, 9).Value = 2434.7727  # BSM!I134: 2528.3 -> 2434.7727
$ws.Cells.Item(134, 11).Value = 7304.3181  # BSM!K134: 7584.900000000001 -> 7304.3181
$ws.Cells.Item(134, 13).Value = -4769.3181  # BSM!M134: -5049.900000000001 -> -4769.3181

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 15262  # CRP!H28: 15660.5 -> 15262
$ws.Cells.Item(28, 10).Value = 15262  # CRP!J28: 15660.5 -> 15262
$ws.Cells.Item(28, 12).Value = 15262  # CRP!L28: 15660.5 -> 15262
$ws.Cells.Item(28, 14).Value = -15752  # CRP!N28: -16150.5 -> -15752

$ws.Cells.Item(31, 8).Value = 1381.8695  # CRP!H31: 1425.4762 -> 1381.8695
$ws.Cells.Item(31, 9).Value = 1409.6471  # CRP!I31: 1444.75 -> 1409.6471
$ws.Cells.Item(31, 10).Value = 1303.1666  # CRP!J31: 1363.8 -> 1303.1666
$ws.Cells.Item(31, 11).Value = 1409.6471  # CRP!K31: 1444.75 -> 1409.6471
$ws.Cells.Item(31, 12).Value = 1303.1666  # CRP!L31: 1363.8 -> 1303.1666
$ws.Cells.Item(31, 13).Value = -1114.6471  # CRP!M31: -1149.75 -> -1114.6471
$ws.Cells.Item(31, 14).Value = -1893.1666  # CRP!N31: -1953.8 -> -1893.1666

$ws.Cells.Item(32, 8).Value = 2106.6667  # CRP!H32: 1886 -> 2106.6667
$ws.Cells.Item(32, 9).Value = 2106.6667  # CRP!I32: 1607.5 -> 2106.6667
$ws.Cells.Item(32, 10).Value = 0  # CRP!J32: 3000 -> 0
$ws.Cells.Item(32, 11).Value = 2106.6667  # CRP!K32: 1607.5 -> 2106.6667
$ws.Cells.Item(32, 12).Value = 0  # CRP!L32: 3000 -> 0
$ws.Cells.Item(32, 13).Value = -1790.6667  # CRP!M32: -1291.5 -> -1790.6667
$ws.Cells.Item(32, 14).ClearContents()  # CRP!N32: delete (was -3632)

$ws.Cells.Item(34, 8).Value = 1381.8695  # CRP!H34: 1425.4762 -> 1381.8695
$ws.Cells.Item(34, 9).Value = 1409.6471  # CRP!I34: 1444.75 -> 1409.6471
$ws.Cells.Item(34, 10).Value = 1303.1666  # CRP!J34: 1363.8 -> 1303.1666
$ws.Cells.Item(34, 11).Value = 1409.6471  # CRP!K34: 1444.75 -> 1409.6471
$ws.Cells.Item(34, 12).Value = 1303.1666  # CRP!L34: 1363.8 -> 1303.1666
$ws.Cells.Item(34, 13).Value = -1207.6471  # CRP!M34: -1242.75 -> -1207.6471
$ws.Cells.Item(34, 14).Value = -1707.1666  # CRP!N34: -1767.8 -> -1707.1666

$ws.Cells.Item(62, 8).Value = 3743.7368  # CRP!H62: 3543.5293 -> 3743.7368
$ws.Cells.Item(62, 9).Value = 3550.6155  # CRP!I62: 3221.5 -> 3550.6155
$ws.Cells.Item(62, 10).Value = 4162.1665  # CRP!J62: 4316.4 -> 4162.1665
$ws.Cells.Item(62, 11).Value = 3550.6155  # CRP!K62: 3221.5 -> 3550.6155
$ws.Cells.Item(62, 12).Value = 4162.1665  # CRP!L62: 4316.4 -> 4162.1665
$ws.Cells.Item(62, 13).Value = -2926.6155  # CRP!M62: -2597.5 -> -2926.6155
$ws.Cells.Item(62, 14).Value = -5410.1665  # CRP!N62: -5564.4 -> -5410.1665

$ws.Cells.Item(65, 8).Value = 3743.7368  # CRP!H65: 3543.5293 -> 3743.7368
$ws.Cells.Item(65, 9).Value = 3550.6155  # CRP!I65: 3221.5 -> 3550.6155
$ws.Cells.Item(65, 10).Value = 4162.1665  # CRP!J65: 4316.4 -> 4162.1665
$ws.Cells.Item(65, 11).Value = 17753.0775  # CRP!K65: 16107.5 -> 17753.0775
$ws.Cells.Item(65, 12).Value = 20810.8325  # CRP!L65: 21582 -> 20810.8325
$ws.Cells.Item(65, 13).Value = -14633.0775  # CRP!M65: -12987.5 -> -14633.0775
$ws.Cells.Item(65, 14).Value = -27050.8325  # CRP!N65: -27822 -> -27050.8325

$ws.Cells.Item(99, 8).Value = 5308.9165  # CRP!H99: 5746 -> 5308.9165
$ws.Cells.Item(99, 9).Value = 5263.5  # CRP!I99: 5275.875 -> 5263.5
$ws.Cells.Item(99, 10).Value = 5399.75  # CRP!J99: 6999.6665 -> 5399.75
$ws.Cells.Item(99, 11).Value = 5263.5  # CRP!K99: 5275.875 -> 5263.5
$ws.Cells.Item(99, 12).Value = 5399.75  # CRP!L99: 6999.6665 -> 5399.75
$ws.Cells.Item(99, 13).Value = -3765.5  # CRP!M99: -3777.875 -> -3765.5
$ws.Cells.Item(99, 14).Value = -8395.75  # CRP!N99: -9995.666499999999 -> -8395.75

$ws.Cells.Item(105, 8).Value = 1427.75  # CRP!H105: 1402.6552 -> 1427.75
$ws.Cells.Item(105, 9).Value = 1476  # CRP!I105: 1447.2593 -> 1476
$ws.Cells.Item(105, 11).Value = 1476  # CRP!K105: 1447.2593 -> 1476
$ws.Cells.Item(105, 13).Value = 271  # CRP!M105: 299.7407000000001 -> 271

$ws.Cells.Item(107, 8).Value = 3999.1177  # CRP!H107: 3007.5417 -> 3999.1177
$ws.Cells.Item(107, 9).Value = 0  # CRP!I107: 609.4 -> 0
$ws.Cells.Item(107, 10).Value = 3999.1177  # CRP!J107: 3638.6316 -> 3999.1177
$ws.Cells.Item(107, 11).Value = 0  # CRP!K107: 609.4 -> 0
$ws.Cells.Item(107, 12).Value = 3999.1177  # CRP!L107: 3638.6316 -> 3999.1177
$ws.Cells.Item(107, 13).ClearContents()  # CRP!M107: delete (was 1310.6)
$ws.Cells.Item(107, 14).Value = -7839.1177  # CRP!N107: -7478.631600000001 -> -7839.1177

$ws.Cells.Item(126, 8).Value = 5308.9165  # CRP!H126: 5746 -> 5308.9165
$ws.Cells.Item(126, 9).Value = 5263.5  # CRP!I126: 5275.875 -> 5263.5
$ws.Cells.Item(126, 10).Value = 5399.75  # CRP!J126: 6999.6665 -> 5399.75
$ws.Cells.Item(126, 11).Value = 15790.5  # CRP!K126: 15827.625 -> 15790.5
$ws.Cells.Item(126, 12).Value = 16199.25  # CRP!L126: 20998.9995 -> 16199.25
$ws.Cells.Item(126, 13).Value = -13320.5  # CRP!M126: -13357.625 -> -13320.5
$ws.Cells.Item(126, 14).Value = -21139.25  # CRP!N126: -25938.9995 -> -21139.25

$ws.Cells.Item(134, 8).Value = 52021.45  # CRP!H134: 57690.723 -> 52021.45
$ws.Cells.Item(134, 9).Value = 78579.38  # CRP!I134: 92685.09 -> 78579.38
$ws.Cells.Item(134, 11).Value = 235738.14  # CRP!K134: 278055.27 -> 235738.14
$ws.Cells.Item(134, 13).Value = -233203.14  # CRP!M134: -275520.27 -> -233203.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 412  # CUL!H11: 412.33334 -> 412
$ws.Cells.Item(11, 9).Value = 412  # CUL!I11: 412.33334 -> 412
$ws.Cells.Item(11, 11).Value = 1236  # CUL!K11: 1237.00002 -> 1236
$ws.Cells.Item(11, 13).Value = -1096  # CUL!M11: -1097.00002 -> -1096

$ws.Cells.Item(33, 8).Value = 67  # CUL!H33: 84.5 -> 67
$ws.Cells.Item(33, 9).Value = 65.454544  # CUL!I33: 85 -> 65.454544
$ws.Cells.Item(33, 11).Value = 392.727264  # CUL!K33: 510 -> 392.727264
$ws.Cells.Item(33, 13).Value = -109.727264  # CUL!M33: -227 -> -109.727264

$ws.Cells.Item(38, 8).Value = 514.931  # CUL!H38: 501.03333 -> 514.931
$ws.Cells.Item(38, 10).Value = 523  # CUL!J38: 462.2857 -> 523
$ws.Cells.Item(38, 12).Value = 1569  # CUL!L38: 1386.8571 -> 1569
$ws.Cells.Item(38, 14).Value = -2263  # CUL!N38: -2080.8571 -> -2263

$ws.Cells.Item(40, 8).Value = 233.33333  # CUL!H40: 131.125 -> 233.33333
$ws.Cells.Item(40, 9).Value = 0  # CUL!I40: 37.5 -> 0
$ws.Cells.Item(40, 10).Value = 233.33333  # CUL!J40: 224.75 -> 233.33333
$ws.Cells.Item(40, 11).Value = 0  # CUL!K40: 150 -> 0
$ws.Cells.Item(40, 12).Value = 933.33332  # CUL!L40: 899 -> 933.33332
$ws.Cells.Item(40, 13).ClearContents()  # CUL!M40: delete (was -81)
$ws.Cells.Item(40, 14).Value = -1071.33332  # CUL!N40: -1037 -> -1071.33332

$ws.Cells.Item(68, 8).Value = 2949  # CUL!H68: 2619.2 -> 2949
$ws.Cells.Item(68, 10).Value = 2949  # CUL!J68: 2619.2 -> 2949
$ws.Cells.Item(68, 12).Value = 8847  # CUL!L68: 7857.599999999999 -> 8847
$ws.Cells.Item(68, 14).Value = -10469  # CUL!N68: -9479.599999999999 -> -10469

$ws.Cells.Item(71, 8).Value = 2949  # CUL!H71: 2619.2 -> 2949
$ws.Cells.Item(71, 10).Value = 2949  # CUL!J71: 2619.2 -> 2949
$ws.Cells.Item(71, 12).Value = 26541  # CUL!L71: 23572.8 -> 26541
$ws.Cells.Item(71, 14).Value = -34653  # CUL!N71: -31684.8 -> -34653

$ws.Cells.Item(80, 8).Value = 85394.8  # CUL!H80: 85396.39999999999 -> 85394.8
$ws.Cells.Item(80, 9).Value = 8993  # CUL!I80: 8993.5 -> 8993
$ws.Cells.Item(80, 10).Value = 199997.5  # CUL!J80: 136331.67 -> 199997.5
$ws.Cells.Item(80, 11).Value = 26979  # CUL!K80: 26980.5 -> 26979
$ws.Cells.Item(80, 12).Value = 599992.5  # CUL!L80: 408995.01 -> 599992.5
$ws.Cells.Item(80, 13).Value = -26043  # CUL!M80: -26044.5 -> -26043
$ws.Cells.Item(80, 14).Value = -601864.5  # CUL!N80: -410867.01 -> -601864.5

$ws.Cells.Item(83, 8).Value = 85394.8  # CUL!H83: 85396.39999999999 -> 85394.8
$ws.Cells.Item(83, 9).Value = 8993  # CUL!I83: 8993.5 -> 8993
$ws.Cells.Item(83, 10).Value = 199997.5  # CUL!J83: 136331.67 -> 199997.5
$ws.Cells.Item(83, 11).Value = 80937  # CUL!K83: 80941.5 -> 80937
$ws.Cells.Item(83, 12).Value = 1799977.5  # CUL!L83: 1226985.03 -> 1799977.5
$ws.Cells.Item(83, 13).Value = -76257  # CUL!M83: -76261.5 -> -76257
$ws.Cells.Item(83, 14).Value = -1809337.5  # CUL!N83: -1236345.03 -> -1809337.5

$ws.Cells.Item(107, 8).Value = 297.5  # CUL!H107: 275.44446 -> 297.5
$ws.Cells.Item(107, 10).Value = 495  # CUL!J107: 415.8 -> 495
$ws.Cells.Item(107, 12).Value = 1485  # CUL!L107: 1247.4 -> 1485
$ws.Cells.Item(107, 14).Value = -5325  # CUL!N107: -5087.4 -> -5325

$ws.Cells.Item(122, 8).Value = 686.7857  # CUL!H122: 710.2308 -> 686.7857
$ws.Cells.Item(122, 9).Value = 450.85715  # CUL!I122: 462.33334 -> 450.85715
$ws.Cells.Item(122, 11).Value = 4057.71435  # CUL!K122: 4161.00006 -> 4057.71435
$ws.Cells.Item(122, 13).Value = -1607.71435  # CUL!M122: -1711.00006 -> -1607.71435

$ws.Cells.Item(131, 8).Value = 16665.715  # CUL!H131: 23224.555 -> 16665.715
$ws.Cells.Item(131, 10).Value = 22929.9  # CUL!J131: 41000 -> 22929.9
$ws.Cells.Item(131, 12).Value = 68789.70000000001  # CUL!L131: 123000 -> 68789.70000000001
$ws.Cells.Item(131, 14).Value = -78869.70000000001  # CUL!N131: -133080 -> -78869.70000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8499.5  # GSM!H70: 6874 -> 8499.5
$ws.Cells.Item(70, 9).Value = 0  # GSM!I70: 4997 -> 0
$ws.Cells.Item(70, 10).Value = 8499.5  # GSM!J70: 7499.6665 -> 8499.5
$ws.Cells.Item(70, 11).Value = 0  # GSM!K70: 4997 -> 0
$ws.Cells.Item(70, 12).Value = 8499.5  # GSM!L70: 7499.6665 -> 8499.5
$ws.Cells.Item(70, 13).ClearContents()  # GSM!M70: delete (was -4727)
$ws.Cells.Item(70, 14).Value = -9039.5  # GSM!N70: -8039.6665 -> -9039.5

$ws.Cells.Item(73, 8).Value = 8499.5  # GSM!H73: 6874 -> 8499.5
$ws.Cells.Item(73, 9).Value = 0  # GSM!I73: 4997 -> 0
$ws.Cells.Item(73, 10).Value = 8499.5  # GSM!J73: 7499.6665 -> 8499.5
$ws.Cells.Item(73, 11).Value = 0  # GSM!K73: 4997 -> 0
$ws.Cells.Item(73, 12).Value = 8499.5  # GSM!L73: 7499.6665 -> 8499.5
$ws.Cells.Item(73, 13).ClearContents()  # GSM!M73: delete (was -4061)
$ws.Cells.Item(73, 14).Value = -10371.5  # GSM!N73: -9371.666499999999 -> -10371.5

$ws.Cells.Item(80, 8).Value = 3761.5  # GSM!H80: 3857.5715 -> 3761.5
$ws.Cells.Item(80, 9).Value = 3218.4  # GSM!I80: 3250.75 -> 3218.4
$ws.Cells.Item(80, 11).Value = 3218.4  # GSM!K80: 3250.75 -> 3218.4
$ws.Cells.Item(80, 13).Value = -2220.4  # GSM!M80: -2252.75 -> -2220.4

$ws.Cells.Item(83, 8).Value = 3761.5  # GSM!H83: 3857.5715 -> 3761.5
$ws.Cells.Item(83, 9).Value = 3218.4  # GSM!I83: 3250.75 -> 3218.4
$ws.Cells.Item(83, 11).Value = 16092  # GSM!K83: 16253.75 -> 16092
$ws.Cells.Item(83, 13).Value = -11100  # GSM!M83: -11261.75 -> -11100

$ws.Cells.Item(97, 8).Value = 1174  # GSM!H97: 1174.875 -> 1174
$ws.Cells.Item(97, 9).Value = 921.13635  # GSM!I97: 922.0909 -> 921.13635
$ws.Cells.Item(97, 11).Value = 921.13635  # GSM!K97: 922.0909 -> 921.13635
$ws.Cells.Item(97, 13).Value = -425.13635  # GSM!M97: -426.0909 -> -425.13635

$ws.Cells.Item(102, 8).Value = 4241.222  # GSM!H102: 4508.353 -> 4241.222
$ws.Cells.Item(102, 9).Value = 3212.3333  # GSM!I102: 3531.6365 -> 3212.3333
$ws.Cells.Item(102, 11).Value = 3212.3333  # GSM!K102: 3531.6365 -> 3212.3333
$ws.Cells.Item(102, 13).Value = -1590.3333  # GSM!M102: -1909.6365 -> -1590.3333

$ws.Cells.Item(113, 8).Value = 120969.766  # GSM!H113: 121065.35 -> 120969.766
$ws.Cells.Item(113, 9).Value = 102771.6  # GSM!I113: 93884.27 -> 102771.6
$ws.Cells.Item(113, 10).Value = 146967.14  # GSM!J113: 170897.33 -> 146967.14
$ws.Cells.Item(113, 11).Value = 102771.6  # GSM!K113: 93884.27 -> 102771.6
$ws.Cells.Item(113, 12).Value = 146967.14  # GSM!L113: 170897.33 -> 146967.14
$ws.Cells.Item(113, 13).Value = -100601.6  # GSM!M113: -91714.27 -> -100601.6
$ws.Cells.Item(113, 14).Value = -151307.14  # GSM!N113: -175237.33 -> -151307.14

$ws.Cells.Item(122, 8).Value = 3624.5217  # GSM!H122: 3770.238 -> 3624.5217
$ws.Cells.Item(122, 10).Value = 4778.8335  # GSM!J122: 5315.7 -> 4778.8335
$ws.Cells.Item(122, 12).Value = 14336.5005  # GSM!L122: 15947.1 -> 14336.5005
$ws.Cells.Item(122, 14).Value = -19236.5005  # GSM!N122: -20847.1 -> -19236.5005

$ws.Cells.Item(132, 8).Value = 203179  # GSM!H132: 203169 -> 203179
$ws.Cells.Item(132, 9).Value = 203179  # GSM!I132: 203169 -> 203179
$ws.Cells.Item(132, 11).Value = 609537  # GSM!K132: 609507 -> 609537
$ws.Cells.Item(132, 13).Value = -607007  # GSM!M132: -606977 -> -607007

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4907.3076  # LTW!H7: 6742.5713 -> 4907.3076
$ws.Cells.Item(7, 9).Value = 3254.0908  # LTW!I7: 3699.6667 -> 3254.0908
$ws.Cells.Item(7, 10).Value = 14000  # LTW!J7: 25000 -> 14000
$ws.Cells.Item(7, 11).Value = 3254.0908  # LTW!K7: 3699.6667 -> 3254.0908
$ws.Cells.Item(7, 12).Value = 14000  # LTW!L7: 25000 -> 14000
$ws.Cells.Item(7, 13).Value = -3142.0908  # LTW!M7: -3587.6667 -> -3142.0908
$ws.Cells.Item(7, 14).Value = -14224  # LTW!N7: -25224 -> -14224

$ws.Cells.Item(16, 8).Value = 2075.375  # LTW!H16: 2080.625 -> 2075.375
$ws.Cells.Item(16, 9).Value = 2074.8572  # LTW!I16: 2078.0356 -> 2074.8572
$ws.Cells.Item(16, 10).Value = 2079  # LTW!J16: 2098.75 -> 2079
$ws.Cells.Item(16, 11).Value = 2074.8572  # LTW!K16: 2078.0356 -> 2074.8572
$ws.Cells.Item(16, 12).Value = 2079  # LTW!L16: 2098.75 -> 2079
$ws.Cells.Item(16, 13).Value = -1904.8572  # LTW!M16: -1908.0356 -> -1904.8572
$ws.Cells.Item(16, 14).Value = -2419  # LTW!N16: -2438.75 -> -2419

$ws.Cells.Item(40, 8).Value = 9863  # LTW!H40: 12499.857 -> 9863
$ws.Cells.Item(40, 9).Value = 8199  # LTW!I40: 0 -> 8199
$ws.Cells.Item(40, 10).Value = 11249.667  # LTW!J40: 12499.857 -> 11249.667
$ws.Cells.Item(40, 11).Value = 8199  # LTW!K40: 0 -> 8199
$ws.Cells.Item(40, 12).Value = 11249.667  # LTW!L40: 12499.857 -> 11249.667
$ws.Cells.Item(40, 13).Value = -8063  # LTW!M40: None -> -8063
$ws.Cells.Item(40, 14).Value = -11521.667  # LTW!N40: -12771.857 -> -11521.667

$ws.Cells.Item(82, 8).Value = 2433.2173  # LTW!H82: 2399.4736 -> 2433.2173
$ws.Cells.Item(82, 9).Value = 2212.8572  # LTW!I82: 2365.1667 -> 2212.8572
$ws.Cells.Item(82, 10).Value = 2776  # LTW!J82: 2458.2856 -> 2776
$ws.Cells.Item(82, 11).Value = 2212.8572  # LTW!K82: 2365.1667 -> 2212.8572
$ws.Cells.Item(82, 12).Value = 2776  # LTW!L82: 2458.2856 -> 2776
$ws.Cells.Item(82, 13).Value = -1851.8572  # LTW!M82: -2004.1667 -> -1851.8572
$ws.Cells.Item(82, 14).Value = -3498  # LTW!N82: -3180.2856 -> -3498

$ws.Cells.Item(85, 8).Value = 2433.2173  # LTW!H85: 2399.4736 -> 2433.2173
$ws.Cells.Item(85, 9).Value = 2212.8572  # LTW!I85: 2365.1667 -> 2212.8572
$ws.Cells.Item(85, 10).Value = 2776  # LTW!J85: 2458.2856 -> 2776
$ws.Cells.Item(85, 11).Value = 2212.8572  # LTW!K85: 2365.1667 -> 2212.8572
$ws.Cells.Item(85, 12).Value = 2776  # LTW!L85: 2458.2856 -> 2776
$ws.Cells.Item(85, 13).Value = -964.8571999999999  # LTW!M85: -1117.1667 -> -964.8571999999999
$ws.Cells.Item(85, 14).Value = -5272  # LTW!N85: -4954.2856 -> -5272

$ws.Cells.Item(100, 8).Value = 1945.5834  # LTW!H100: 1988.9565 -> 1945.5834
$ws.Cells.Item(100, 10).Value = 3711  # LTW!J100: 4105.7144 -> 3711
$ws.Cells.Item(100, 12).Value = 3711  # LTW!L100: 4105.7144 -> 3711
$ws.Cells.Item(100, 14).Value = -4793  # LTW!N100: -5187.7144 -> -4793

$ws.Cells.Item(122, 8).Value = 4242.8184  # LTW!H122: 4133.6 -> 4242.8184
$ws.Cells.Item(122, 9).Value = 3486.5454  # LTW!I122: 3453.5715 -> 3486.5454
$ws.Cells.Item(122, 11).Value = 10459.6362  # LTW!K122: 10360.7145 -> 10459.6362
$ws.Cells.Item(122, 13).Value = -8009.636200000001  # LTW!M122: -7910.7145 -> -8009.636200000001

$ws.Cells.Item(126, 8).Value = 4907.3076  # LTW!H126: 6742.5713 -> 4907.3076
$ws.Cells.Item(126, 9).Value = 3254.0908  # LTW!I126: 3699.6667 -> 3254.0908
$ws.Cells.Item(126, 10).Value = 14000  # LTW!J126: 25000 -> 14000
$ws.Cells.Item(126, 11).Value = 9762.2724  # LTW!K126: 11099.0001 -> 9762.2724
$ws.Cells.Item(126, 12).Value = 42000  # LTW!L126: 75000 -> 42000
$ws.Cells.Item(126, 13).Value = -7292.2724  # LTW!M126: -8629.000100000001 -> -7292.2724
$ws.Cells.Item(126, 14).Value = -46940  # LTW!N126: -79940 -> -46940

$ws.Cells.Item(132, 8).Value = 97950.84  # LTW!H132: 98049.46000000001 -> 97950.84
$ws.Cells.Item(132, 9).Value = 97950.84  # LTW!I132: 98049.46000000001 -> 97950.84
$ws.Cells.Item(132, 11).Value = 293852.52  # LTW!K132: 294148.38 -> 293852.52
$ws.Cells.Item(132, 13).Value = -291322.52  # LTW!M132: -291618.38 -> -291322.52

$ws.Cells.Item(133, 8).Value = 89978.664  # LTW!H133: 89988 -> 89978.664
$ws.Cells.Item(133, 10).Value = 89978.664  # LTW!J133: 89988 -> 89978.664
$ws.Cells.Item(133, 12).Value = 89978.664  # LTW!L133: 89988 -> 89978.664
$ws.Cells.Item(133, 14).Value = -95038.664  # LTW!N133: -95048 -> -95038.664

$ws.Cells.Item(136, 8).Value = 3111.3333  # LTW!H136: 3363.1875 -> 3111.3333
$ws.Cells.Item(136, 9).Value = 2469  # LTW!I136: 2605.1538 -> 2469
$ws.Cells.Item(136, 10).Value = 8250  # LTW!J136: 6648 -> 8250
$ws.Cells.Item(136, 11).Value = 7407  # LTW!K136: 7815.4614 -> 7407
$ws.Cells.Item(136, 12).Value = 24750  # LTW!L136: 19944 -> 24750
$ws.Cells.Item(136, 13).Value = -4857  # LTW!M136: -5265.4614 -> -4857
$ws.Cells.Item(136, 14).Value = -29850  # LTW!N136: -25044 -> -29850

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 256498.5  # WVR!H62: 339664.66 -> 256498.5
$ws.Cells.Item(62, 9).Value = 7000  # WVR!I62: 0 -> 7000
$ws.Cells.Item(62, 11).Value = 7000  # WVR!K62: 0 -> 7000
$ws.Cells.Item(62, 13).Value = -6376  # WVR!M62: None -> -6376

$ws.Cells.Item(65, 8).Value = 256498.5  # WVR!H65: 339664.66 -> 256498.5
$ws.Cells.Item(65, 9).Value = 7000  # WVR!I65: 0 -> 7000
$ws.Cells.Item(65, 11).Value = 35000  # WVR!K65: 0 -> 35000
$ws.Cells.Item(65, 13).Value = -31880  # WVR!M65: None -> -31880

$ws.Cells.Item(122, 8).Value = 1892.3334  # WVR!H122: 2129.375 -> 1892.3334
$ws.Cells.Item(122, 9).Value = 1892.3334  # WVR!I122: 2129.375 -> 1892.3334
$ws.Cells.Item(122, 11).Value = 5677.0002  # WVR!K122: 6388.125 -> 5677.0002
$ws.Cells.Item(122, 13).Value = -3227.0002  # WVR!M122: -3938.125 -> -3227.0002

$ws.Cells.Item(132, 8).Value = 30070.555  # WVR!H132: 31795.383 -> 30070.555
$ws.Cells.Item(132, 9).Value = 30844  # WVR!I132: 32667.969 -> 30844
$ws.Cells.Item(132, 11).Value = 92532  # WVR!K132: 98003.90700000001 -> 92532
$ws.Cells.Item(132, 13).Value = -90002  # WVR!M132: -95473.90700000001 -> -90002

$ws.Cells.Item(136, 8).Value = 5527.696  # WVR!H136: 5754.1816 -> 5527.696
$ws.Cells.Item(136, 9).Value = 5815.4  # WVR!I136: 6191.857 -> 5815.4
$ws.Cells.Item(136, 11).Value = 17446.2  # WVR!K136: 18575.571 -> 17446.2
$ws.Cells.Item(136, 13).Value = -14896.2  # WVR!M136: -16025.571 -> -14896.2

